# "Minor typo fix on slides" -- CSC301 TDSB Education app final presentation
#
# Slide 8 ("Technical Item 3 - Internal APIs"), Content Placeholder 4:
#   "FirestoreManager - DAO that communicates with out database"
#     -> "FirestoreManager - DAO that communicates without database"
#
# While correcting that bullet, the author also split the leading
# spell-checked term out of the "ClassroomApiAccess ..." and
# "FirestoreManager ..." bullets into their own runs, so we reproduce that
# run layout here as well.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item("Content Placeholder 4")
$tr = $shape.TextFrame.TextRange

$enDash = [char]0x2013

# Paragraph 5: "ClassroomApiAccess <dash> Wrapper for the Google Classroom
# API, allows us to pull information from Classroom"
$para5 = $tr.Paragraphs(5)
$run5 = $para5.Runs(1)
$run5.Text = "ClassroomApiAccess"
$null = $run5.InsertAfter(" " + $enDash + " Wrapper for the Google Classroom API, allows us to pull information from Classroom")

# Paragraph 7: "FirestoreManager <dash> DAO that communicates with out
# database" -> fix the "with out" typo to "without"
$para7 = $tr.Paragraphs(7)
$run7 = $para7.Runs(1)
$run7.Text = "FirestoreManager"
$null = $run7.InsertAfter(" " + $enDash + " DAO that communicates without database")

Write-Output ("Paragraph 5: " + $tr.Paragraphs(5).Text)
Write-Output ("Paragraph 7: " + $tr.Paragraphs(7).Text)
